$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 12
$ws.Range("B1").Value = 12
$ws.Range("C1").Value = 13

$ws.Range("C1").Select()
